# Progress & Challenges slide (slide 3): add a new "body" placeholder shape
# with the progress notes, and remove the old empty "Subtitle 2" placeholder.

$p  = $ppt.ActivePresentation
$s3 = $p.Slides.Item(3)

# ------------------------------------------------------------------
# 1. Obtain a genuine "body" placeholder shape. The current slide layout
#    ("Title Slide") only defines title/subTitle placeholders, so we borrow
#    a body placeholder from a sibling "Title, Content" layout by briefly
#    creating a throw-away slide, copying its body placeholder onto slide 3,
#    then discarding the throw-away slide again.
# ------------------------------------------------------------------
$allLayouts = $p.SlideMaster.CustomLayouts
$contentLayout = $null
for ($i = 1; $i -le $allLayouts.Count; $i++) {
    $cand = $allLayouts.Item($i)
    if ($cand.Name -eq "Title, Content") {
        $contentLayout = $cand
        break
    }
}

$tempSlide = $p.Slides.AddSlide($p.Slides.Count + 1, $contentLayout)

$bodyPlaceholder = $null
for ($i = 1; $i -le $tempSlide.Shapes.Count; $i++) {
    $cand = $tempSlide.Shapes.Item($i)
    if ($cand.PlaceholderFormat.Type -eq 2) {
        $bodyPlaceholder = $cand
    }
}

$bodyPlaceholder.Copy()
$newShape = $s3.Shapes.Paste().Item(1)
$tempSlide.Delete()

$newShape.Name = "Text Placeholder 3"
$newShape.ZOrder(1)   # msoSendToBack -> becomes the first shape on the slide

# ------------------------------------------------------------------
# 2. Fill in the body text (7 paragraphs; two of them indented to level 2).
# ------------------------------------------------------------------
$tf = $newShape.TextFrame
$tr = $tf.TextRange

$tr.Text = "Finished converting Matlab code to Python"
$tr.InsertAfter("`rTraffic generation, collision detection, car management")
$tr.InsertAfter("`rChallenges: Code was not robust and needed to be redesigned more than we expected")
$tr.InsertAfter("`rDeveloped simple visualization for testing")
$tr.InsertAfter("`rImplemented Dresner/Stone and stop sign policies.")
$tr.InsertAfter("`rChallenges: More advanced traffic light policy has taken more time to formulate")
$tr.InsertAfter("`r")

$tr.Paragraphs(2, 1).IndentLevel = 2
$tr.Paragraphs(3, 1).IndentLevel = 2
$tr.Paragraphs(6, 1).IndentLevel = 2

# Split runs within paragraphs that mix plain and "flagged" text (re-assigning
# the same characters forces PowerPoint to materialize separate runs).
$p1 = $tr.Paragraphs(1, 1)
$p1.Characters(1, 20).Text  = "Finished converting "
$p1.Characters(21, 6).Text  = "Matlab"
$p1.Characters(27, 15).Text = " code to Python"

$p5 = $tr.Paragraphs(5, 1)
$p5.Characters(1, 12).Text  = "Implemented "
$p5.Characters(13, 7).Text  = "Dresner"
$p5.Characters(20, 31).Text = "/Stone and stop sign policies."

Write-Host "Final shapes on slide 3:"
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    Write-Host $i $s3.Shapes.Item($i).Name
}

# ------------------------------------------------------------------
# 3. Remove the old, empty "Subtitle 2" placeholder shape.
# ------------------------------------------------------------------
for ($i = $s3.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s3.Shapes.Item($i)
    if ($sh.Name -eq "Subtitle 2") {
        $sh.Delete()
    }
}

Write-Host "Shapes after removing Subtitle 2:"
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    Write-Host $i $s3.Shapes.Item($i).Name
}
